# The workbook has two logical sheets: "Sheet2" (tab order 0, the front
# "Vendas" summary sheet) and "Sheet1" (tab order 1, the detailed
# day-by-day sales data sheet with columns 1-31 and a Total column).
# The diff's edits (column AA, rows 2,4,5,6,7,11,12,13) belong to the
# detailed data sheet, which is the logical worksheet named "Sheet1".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AA2").Value = 32
$ws.Range("AA4").Value = 12
$ws.Range("AA5").Value = 63
$ws.Range("AA6").Value = 27
$ws.Range("AA7").Value = 52
$ws.Range("AA11").Value = 27
$ws.Range("AA12").Value = 12
$ws.Range("AA13").Value = 60
